$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 428
$ws.Range("A428").Value = 10
$ws.Range("B428").Value = "Vega Modelo de Temuco"
$ws.Range("C428").Value = "La Araucanía"
$ws.Range("D428").Value = 44714
$ws.Range("E428").Value = 9
$ws.Range("F428").Value = 100112032
$ws.Range("G428").Value = "Zapallo italiano"
$ws.Range("H428").Value = "Sin especificar"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 165
$ws.Range("K428").Value = 18000
$ws.Range("L428").Value = 18000
$ws.Range("M428").Value = 18000
$ws.Range("N428").Value = "`$/caja 60 unidades"
$ws.Range("O428").Value = "Región de Arica y Parinacota"
$ws.Range("P428").Value = 300
$ws.Range("Q428").Value = 60
$ws.Range("R428").Value = "Hortaliza"

# Row 429
$ws.Range("A429").Value = 10
$ws.Range("B429").Value = "Vega Modelo de Temuco"
$ws.Range("C429").Value = "La Araucanía"
$ws.Range("D429").Value = 44714
$ws.Range("E429").Value = 9
$ws.Range("F429").Value = 100112032
$ws.Range("G429").Value = "Zapallo italiano"
$ws.Range("H429").Value = "Sin especificar"
$ws.Range("I429").Value = "Segunda"
$ws.Range("J429").Value = 95
$ws.Range("K429").Value = 16000
$ws.Range("L429").Value = 16000
$ws.Range("M429").Value = 16000
$ws.Range("N429").Value = "`$/caja 80 unidades"
$ws.Range("O429").Value = "Región de Arica y Parinacota"
$ws.Range("P429").Value = 200
$ws.Range("Q429").Value = 80
$ws.Range("R429").Value = "Hortaliza"

# Row 430
$ws.Range("A430").Value = 10
$ws.Range("B430").Value = "Vega Modelo de Temuco"
$ws.Range("C430").Value = "La Araucanía"
$ws.Range("D430").Value = 44188
$ws.Range("E430").Value = 9
$ws.Range("F430").Value = 100112032
$ws.Range("G430").Value = "Zapallo italiano"
$ws.Range("H430").Value = "Sin especificar"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 400
$ws.Range("K430").Value = 10000
$ws.Range("L430").Value = 10000
$ws.Range("M430").Value = 10000
$ws.Range("N430").Value = "`$/caja 60 unidades"
$ws.Range("O430").Value = "Región del Maule"
$ws.Range("P430").Value = 167
$ws.Range("Q430").Value = 60
$ws.Range("R430").Value = "Hortaliza"

# Row 431
$ws.Range("A431").Value = 10
$ws.Range("B431").Value = "Vega Modelo de Temuco"
$ws.Range("C431").Value = "La Araucanía"
$ws.Range("D431").Value = 44659
$ws.Range("E431").Value = 9
$ws.Range("F431").Value = 100112032
$ws.Range("G431").Value = "Zapallo italiano"
$ws.Range("H431").Value = "Bola 8"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 75
$ws.Range("K431").Value = 15000
$ws.Range("L431").Value = 15000
$ws.Range("M431").Value = 15000
$ws.Range("N431").Value = "`$/caja 60 unidades"
$ws.Range("O431").Value = "Región Metropolitana"
$ws.Range("P431").Value = 250
$ws.Range("Q431").Value = 60
$ws.Range("R431").Value = "Hortaliza"

# Row 432
$ws.Range("A432").Value = 10
$ws.Range("B432").Value = "Vega Modelo de Temuco"
$ws.Range("C432").Value = "La Araucanía"
$ws.Range("D432").Value = 44659
$ws.Range("E432").Value = 9
$ws.Range("F432").Value = 100112032
$ws.Range("G432").Value = "Zapallo italiano"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 125
$ws.Range("K432").Value = 10000
$ws.Range("L432").Value = 10000
$ws.Range("M432").Value = 10000
$ws.Range("N432").Value = "`$/caja 60 unidades"
$ws.Range("O432").Value = "Región del Maule"
$ws.Range("P432").Value = 167
$ws.Range("Q432").Value = 60
$ws.Range("R432").Value = "Hortaliza"

# Row 433
$ws.Range("A433").Value = 10
$ws.Range("B433").Value = "Vega Modelo de Temuco"
$ws.Range("C433").Value = "La Araucanía"
$ws.Range("D433").Value = 44494
$ws.Range("E433").Value = 9
$ws.Range("F433").Value = 100112032
$ws.Range("G433").Value = "Zapallo italiano"
$ws.Range("H433").Value = "Bola 8"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 50
$ws.Range("K433").Value = 15000
$ws.Range("L433").Value = 15000
$ws.Range("M433").Value = 15000
$ws.Range("N433").Value = "`$/caja 60 unidades"
$ws.Range("O433").Value = "Región de O'Higgins"
$ws.Range("P433").Value = 250
$ws.Range("Q433").Value = 60
$ws.Range("R433").Value = "Hortaliza"

# Row 434
$ws.Range("A434").Value = 10
$ws.Range("B434").Value = "Vega Modelo de Temuco"
$ws.Range("C434").Value = "La Araucanía"
$ws.Range("D434").Value = 44494
$ws.Range("E434").Value = 9
$ws.Range("F434").Value = 100112032
$ws.Range("G434").Value = "Zapallo italiano"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 100
$ws.Range("K434").Value = 10000
$ws.Range("L434").Value = 10000
$ws.Range("M434").Value = 10000
$ws.Range("N434").Value = "`$/caja 60 unidades"
$ws.Range("O434").Value = "Región de Arica y Parinacota"
$ws.Range("P434").Value = 167
$ws.Range("Q434").Value = 60
$ws.Range("R434").Value = "Hortaliza"

# Row 435
$ws.Range("A435").Value = 10
$ws.Range("B435").Value = "Vega Modelo de Temuco"
$ws.Range("C435").Value = "La Araucanía"
$ws.Range("D435").Value = 44494
$ws.Range("E435").Value = 9
$ws.Range("F435").Value = 100112032
$ws.Range("G435").Value = "Zapallo italiano"
$ws.Range("H435").Value = "Sin especificar"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 200
$ws.Range("K435").Value = 15000
$ws.Range("L435").Value = 15000
$ws.Range("M435").Value = 15000
$ws.Range("N435").Value = "`$/caja 60 unidades"
$ws.Range("O435").Value = "Región de O'Higgins"
$ws.Range("P435").Value = 250
$ws.Range("Q435").Value = 60
$ws.Range("R435").Value = "Hortaliza"

# Row 436
$ws.Range("A436").Value = 10
$ws.Range("B436").Value = "Vega Modelo de Temuco"
$ws.Range("C436").Value = "La Araucanía"
$ws.Range("D436").Value = 44453
$ws.Range("E436").Value = 9
$ws.Range("F436").Value = 100112032
$ws.Range("G436").Value = "Zapallo italiano"
$ws.Range("H436").Value = "Sin especificar"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 90
$ws.Range("K436").Value = 18000
$ws.Range("L436").Value = 19000
$ws.Range("M436").Value = 18444
$ws.Range("N436").Value = "`$/caja 60 unidades"
$ws.Range("O436").Value = "Región de Arica y Parinacota"
$ws.Range("P436").Value = 307
$ws.Range("Q436").Value = 60
$ws.Range("R436").Value = "Hortaliza"

# Row 437
$ws.Range("A437").Value = 10
$ws.Range("B437").Value = "Vega Modelo de Temuco"
$ws.Range("C437").Value = "La Araucanía"
$ws.Range("D437").Value = 44421
$ws.Range("E437").Value = 9
$ws.Range("F437").Value = 100112032
$ws.Range("G437").Value = "Zapallo italiano"
$ws.Range("H437").Value = "Sin especificar"
$ws.Range("I437").Value = "Primera"
$ws.Range("J437").Value = 200
$ws.Range("K437").Value = 12000
$ws.Range("L437").Value = 12000
$ws.Range("M437").Value = 12000
$ws.Range("N437").Value = "`$/caja 60 unidades"
$ws.Range("O437").Value = "Región de Arica y Parinacota"
$ws.Range("P437").Value = 200
$ws.Range("Q437").Value = 60
$ws.Range("R437").Value = "Hortaliza"

# Row 438
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 44291
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = 100112032
$ws.Range("G438").Value = "Zapallo italiano"
$ws.Range("H438").Value = "Bola 8"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 45
$ws.Range("K438").Value = 12000
$ws.Range("L438").Value = 12000
$ws.Range("M438").Value = 12000
$ws.Range("N438").Value = "`$/caja 60 unidades"
$ws.Range("O438").Value = "Región del Maule"
$ws.Range("P438").Value = 200
$ws.Range("Q438").Value = 60
$ws.Range("R438").Value = "Hortaliza"

# Row 439
$ws.Range("A439").Value = 10
$ws.Range("B439").Value = "Vega Modelo de Temuco"
$ws.Range("C439").Value = "La Araucanía"
$ws.Range("D439").Value = 44291
$ws.Range("E439").Value = 9
$ws.Range("F439").Value = 100112032
$ws.Range("G439").Value = "Zapallo italiano"
$ws.Range("H439").Value = "Sin especificar"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 155
$ws.Range("K439").Value = 13000
$ws.Range("L439").Value = 13000
$ws.Range("M439").Value = 13000
$ws.Range("N439").Value = "`$/caja 60 unidades"
$ws.Range("O439").Value = "Región del Maule"
$ws.Range("P439").Value = 217
$ws.Range("Q439").Value = 60
$ws.Range("R439").Value = "Hortaliza"

# Row 440
$ws.Range("A440").Value = 10
$ws.Range("B440").Value = "Vega Modelo de Temuco"
$ws.Range("C440").Value = "La Araucanía"
$ws.Range("D440").Value = 44323
$ws.Range("E440").Value = 9
$ws.Range("F440").Value = 100112032
$ws.Range("G440").Value = "Zapallo italiano"
$ws.Range("H440").Value = "Bola 8"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 55
$ws.Range("K440").Value = 12000
$ws.Range("L440").Value = 12000
$ws.Range("M440").Value = 12000
$ws.Range("N440").Value = "`$/caja 60 unidades"
$ws.Range("O440").Value = "Región de Arica y Parinacota"
$ws.Range("P440").Value = 200
$ws.Range("Q440").Value = 60
$ws.Range("R440").Value = "Hortaliza"

# Row 441
$ws.Range("A441").Value = 10
$ws.Range("B441").Value = "Vega Modelo de Temuco"
$ws.Range("C441").Value = "La Araucanía"
$ws.Range("D441").Value = 44323
$ws.Range("E441").Value = 9
$ws.Range("F441").Value = 100112032
$ws.Range("G441").Value = "Zapallo italiano"
$ws.Range("H441").Value = "Sin especificar"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 125
$ws.Range("K441").Value = 10000
$ws.Range("L441").Value = 10000
$ws.Range("M441").Value = 10000
$ws.Range("N441").Value = "`$/caja 60 unidades"
$ws.Range("O441").Value = "Región de Arica y Parinacota"
$ws.Range("P441").Value = 167
$ws.Range("Q441").Value = 60
$ws.Range("R441").Value = "Hortaliza"

# Row 442
$ws.Range("A442").Value = 10
$ws.Range("B442").Value = "Vega Modelo de Temuco"
$ws.Range("C442").Value = "La Araucanía"
$ws.Range("D442").Value = 44526
$ws.Range("E442").Value = 9
$ws.Range("F442").Value = 100112032
$ws.Range("G442").Value = "Zapallo italiano"
$ws.Range("H442").Value = "Bola 8"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 40
$ws.Range("K442").Value = 8000
$ws.Range("L442").Value = 8000
$ws.Range("M442").Value = 8000
$ws.Range("N442").Value = "`$/caja 60 unidades"
$ws.Range("O442").Value = "Región de Arica y Parinacota"
$ws.Range("P442").Value = 133
$ws.Range("Q442").Value = 60
$ws.Range("R442").Value = "Hortaliza"

# Row 443
$ws.Range("A443").Value = 10
$ws.Range("B443").Value = "Vega Modelo de Temuco"
$ws.Range("C443").Value = "La Araucanía"
$ws.Range("D443").Value = 44526
$ws.Range("E443").Value = 9
$ws.Range("F443").Value = 100112032
$ws.Range("G443").Value = "Zapallo italiano"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 290
$ws.Range("K443").Value = 8000
$ws.Range("L443").Value = 9000
$ws.Range("M443").Value = 8483
$ws.Range("N443").Value = "`$/caja 60 unidades"
$ws.Range("O443").Value = "Región del Maule"
$ws.Range("P443").Value = 141
$ws.Range("Q443").Value = 60
$ws.Range("R443").Value = "Hortaliza"

# Row 444
$ws.Range("A444").Value = 10
$ws.Range("B444").Value = "Vega Modelo de Temuco"
$ws.Range("C444").Value = "La Araucanía"
$ws.Range("D444").Value = 44363
$ws.Range("E444").Value = 9
$ws.Range("F444").Value = 100112032
$ws.Range("G444").Value = "Zapallo italiano"
$ws.Range("H444").Value = "Sin especificar"
$ws.Range("I444").Value = "Primera"
$ws.Range("J444").Value = 100
$ws.Range("K444").Value = 12000
$ws.Range("L444").Value = 12000
$ws.Range("M444").Value = 12000
$ws.Range("N444").Value = "`$/caja 60 unidades"
$ws.Range("O444").Value = "Región de Arica y Parinacota"
$ws.Range("P444").Value = 200
$ws.Range("Q444").Value = 60
$ws.Range("R444").Value = "Hortaliza"

# Row 445
$ws.Range("A445").Value = 10
$ws.Range("B445").Value = "Vega Modelo de Temuco"
$ws.Range("C445").Value = "La Araucanía"
$ws.Range("D445").Value = 44251
$ws.Range("E445").Value = 9
$ws.Range("F445").Value = 100112032
$ws.Range("G445").Value = "Zapallo italiano"
$ws.Range("H445").Value = "Sin especificar"
$ws.Range("I445").Value = "Primera"
$ws.Range("J445").Value = 180
$ws.Range("K445").Value = 1000
$ws.Range("L445").Value = 9000
$ws.Range("M445").Value = 5222
$ws.Range("N445").Value = "`$/caja 60 unidades"
$ws.Range("O445").Value = "Región del Maule"
$ws.Range("P445").Value = 87
$ws.Range("Q445").Value = 60
$ws.Range("R445").Value = "Hortaliza"

# Row 446
$ws.Range("A446").Value = 10
$ws.Range("B446").Value = "Vega Modelo de Temuco"
$ws.Range("C446").Value = "La Araucanía"
$ws.Range("D446").Value = 44403
$ws.Range("E446").Value = 9
$ws.Range("F446").Value = 100112032
$ws.Range("G446").Value = "Zapallo italiano"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 300
$ws.Range("K446").Value = 10000
$ws.Range("L446").Value = 12000
$ws.Range("M446").Value = 11000
$ws.Range("N446").Value = "`$/caja 60 unidades"
$ws.Range("O446").Value = "Región de Arica y Parinacota"
$ws.Range("P446").Value = 183
$ws.Range("Q446").Value = 60
$ws.Range("R446").Value = "Hortaliza"

# Row 447
$ws.Range("A447").Value = 10
$ws.Range("B447").Value = "Vega Modelo de Temuco"
$ws.Range("C447").Value = "La Araucanía"
$ws.Range("D447").Value = 44704
$ws.Range("E447").Value = 9
$ws.Range("F447").Value = 100112032
$ws.Range("G447").Value = "Zapallo italiano"
$ws.Range("H447").Value = "Bola 8"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 100
$ws.Range("K447").Value = 17000
$ws.Range("L447").Value = 20000
$ws.Range("M447").Value = 18500
$ws.Range("N447").Value = "`$/caja 60 unidades"
$ws.Range("O447").Value = "Región de Arica y Parinacota"
$ws.Range("P447").Value = 308
$ws.Range("Q447").Value = 60
$ws.Range("R447").Value = "Hortaliza"

# Row 448
$ws.Range("A448").Value = 10
$ws.Range("B448").Value = "Vega Modelo de Temuco"
$ws.Range("C448").Value = "La Araucanía"
$ws.Range("D448").Value = 44704
$ws.Range("E448").Value = 9
$ws.Range("F448").Value = 100112032
$ws.Range("G448").Value = "Zapallo italiano"
$ws.Range("H448").Value = "Sin especificar"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 220
$ws.Range("K448").Value = 19000
$ws.Range("L448").Value = 20000
$ws.Range("M448").Value = 19545
$ws.Range("N448").Value = "`$/caja 60 unidades"
$ws.Range("O448").Value = "Región de Arica y Parinacota"
$ws.Range("P448").Value = 326
$ws.Range("Q448").Value = 60
$ws.Range("R448").Value = "Hortaliza"

# Row 449
$ws.Range("A449").Value = 10
$ws.Range("B449").Value = "Vega Modelo de Temuco"
$ws.Range("C449").Value = "La Araucanía"
$ws.Range("D449").Value = 44620
$ws.Range("E449").Value = 9
$ws.Range("F449").Value = 100112032
$ws.Range("G449").Value = "Zapallo italiano"
$ws.Range("H449").Value = "Sin especificar"
$ws.Range("I449").Value = "Primera"
$ws.Range("J449").Value = 130
$ws.Range("K449").Value = 10000
$ws.Range("L449").Value = 11000
$ws.Range("M449").Value = 10385
$ws.Range("N449").Value = "`$/caja 60 unidades"
$ws.Range("O449").Value = "Región del Maule"
$ws.Range("P449").Value = 173
$ws.Range("Q449").Value = 60
$ws.Range("R449").Value = "Hortaliza"

# Row 450
$ws.Range("A450").Value = 10
$ws.Range("B450").Value = "Vega Modelo de Temuco"
$ws.Range("C450").Value = "La Araucanía"
$ws.Range("D450").Value = 44586
$ws.Range("E450").Value = 9
$ws.Range("F450").Value = 100112032
$ws.Range("G450").Value = "Zapallo italiano"
$ws.Range("H450").Value = "Sin especificar"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 198
$ws.Range("K450").Value = 14000
$ws.Range("L450").Value = 14000
$ws.Range("M450").Value = 14000
$ws.Range("N450").Value = "`$/caja 60 unidades"
$ws.Range("O450").Value = "Región del Maule"
$ws.Range("P450").Value = 233
$ws.Range("Q450").Value = 60
$ws.Range("R450").Value = "Hortaliza"

# Row 451
$ws.Range("A451").Value = 10
$ws.Range("B451").Value = "Vega Modelo de Temuco"
$ws.Range("C451").Value = "La Araucanía"
$ws.Range("D451").Value = 44601
$ws.Range("E451").Value = 9
$ws.Range("F451").Value = 100112032
$ws.Range("G451").Value = "Zapallo italiano"
$ws.Range("H451").Value = "Sin especificar"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 65
$ws.Range("K451").Value = 12000
$ws.Range("L451").Value = 12000
$ws.Range("M451").Value = 12000
$ws.Range("N451").Value = "`$/caja 60 unidades"
$ws.Range("O451").Value = "Región del Maule"
$ws.Range("P451").Value = 200
$ws.Range("Q451").Value = 60
$ws.Range("R451").Value = "Hortaliza"

# Row 452
$ws.Range("A452").Value = 10
$ws.Range("B452").Value = "Vega Modelo de Temuco"
$ws.Range("C452").Value = "La Araucanía"
$ws.Range("D452").Value = 44544
$ws.Range("E452").Value = 9
$ws.Range("F452").Value = 100112032
$ws.Range("G452").Value = "Zapallo italiano"
$ws.Range("H452").Value = "Sin especificar"
$ws.Range("I452").Value = "Primera"
$ws.Range("J452").Value = 223
$ws.Range("K452").Value = 8000
$ws.Range("L452").Value = 9000
$ws.Range("M452").Value = 8439
$ws.Range("N452").Value = "`$/caja 60 unidades"
$ws.Range("O452").Value = "Región del Maule"
$ws.Range("P452").Value = 141
$ws.Range("Q452").Value = 60
$ws.Range("R452").Value = "Hortaliza"

# Row 453
$ws.Range("A453").Value = 10
$ws.Range("B453").Value = "Vega Modelo de Temuco"
$ws.Range("C453").Value = "La Araucanía"
$ws.Range("D453").Value = 44617
$ws.Range("E453").Value = 9
$ws.Range("F453").Value = 100112032
$ws.Range("G453").Value = "Zapallo italiano"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 50
$ws.Range("K453").Value = 10000
$ws.Range("L453").Value = 10000
$ws.Range("M453").Value = 10000
$ws.Range("N453").Value = "`$/caja 60 unidades"
$ws.Range("O453").Value = "Región del Maule"
$ws.Range("P453").Value = 167
$ws.Range("Q453").Value = 60
$ws.Range("R453").Value = "Hortaliza"

# Row 454
$ws.Range("A454").Value = 10
$ws.Range("B454").Value = "Vega Modelo de Temuco"
$ws.Range("C454").Value = "La Araucanía"
$ws.Range("D454").Value = 44567
$ws.Range("E454").Value = 9
$ws.Range("F454").Value = 100112032
$ws.Range("G454").Value = "Zapallo italiano"
$ws.Range("H454").Value = "Sin especificar"
$ws.Range("I454").Value = "Primera"
$ws.Range("J454").Value = 500
$ws.Range("K454").Value = 8000
$ws.Range("L454").Value = 8000
$ws.Range("M454").Value = 8000
$ws.Range("N454").Value = "`$/caja 60 unidades"
$ws.Range("O454").Value = "Región del Maule"
$ws.Range("P454").Value = 133
$ws.Range("Q454").Value = 60
$ws.Range("R454").Value = "Hortaliza"

# Row 455
$ws.Range("A455").Value = 10
$ws.Range("B455").Value = "Vega Modelo de Temuco"
$ws.Range("C455").Value = "La Araucanía"
$ws.Range("D455").Value = 44169
$ws.Range("E455").Value = 9
$ws.Range("F455").Value = 100112032
$ws.Range("G455").Value = "Zapallo italiano"
$ws.Range("H455").Value = "Bola 8"
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 105
$ws.Range("K455").Value = 9000
$ws.Range("L455").Value = 10000
$ws.Range("M455").Value = 9619
$ws.Range("N455").Value = "`$/caja 60 unidades"
$ws.Range("O455").Value = "Región del Maule"
$ws.Range("P455").Value = 160
$ws.Range("Q455").Value = 60
$ws.Range("R455").Value = "Hortaliza"

# Row 456
$ws.Range("A456").Value = 10
$ws.Range("B456").Value = "Vega Modelo de Temuco"
$ws.Range("C456").Value = "La Araucanía"
$ws.Range("D456").Value = 44169
$ws.Range("E456").Value = 9
$ws.Range("F456").Value = 100112032
$ws.Range("G456").Value = "Zapallo italiano"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 155
$ws.Range("K456").Value = 8000
$ws.Range("L456").Value = 8000
$ws.Range("M456").Value = 8000
$ws.Range("N456").Value = "`$/caja 60 unidades"
$ws.Range("O456").Value = "Región del Maule"
$ws.Range("P456").Value = 133
$ws.Range("Q456").Value = 60
$ws.Range("R456").Value = "Hortaliza"

# Row 457
$ws.Range("A457").Value = 10
$ws.Range("B457").Value = "Vega Modelo de Temuco"
$ws.Range("C457").Value = "La Araucanía"
$ws.Range("D457").Value = 44377
$ws.Range("E457").Value = 9
$ws.Range("F457").Value = 100112032
$ws.Range("G457").Value = "Zapallo italiano"
$ws.Range("H457").Value = "Sin especificar"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 110
$ws.Range("K457").Value = 12000
$ws.Range("L457").Value = 12000
$ws.Range("M457").Value = 12000
$ws.Range("N457").Value = "`$/caja 60 unidades"
$ws.Range("O457").Value = "Región de Arica y Parinacota"
$ws.Range("P457").Value = 200
$ws.Range("Q457").Value = 60
$ws.Range("R457").Value = "Hortaliza"

# Row 458
$ws.Range("A458").Value = 10
$ws.Range("B458").Value = "Vega Modelo de Temuco"
$ws.Range("C458").Value = "La Araucanía"
$ws.Range("D458").Value = 44487
$ws.Range("D458").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E458").Value = 9
$ws.Range("F458").Value = 100112032
$ws.Range("G458").Value = "Zapallo italiano"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 155
$ws.Range("K458").Value = 13000
$ws.Range("L458").Value = 13000
$ws.Range("M458").Value = 13000
$ws.Range("N458").Value = "`$/caja 60 unidades"
$ws.Range("O458").Value = "Región de Arica y Parinacota"
$ws.Range("P458").Value = 217
$ws.Range("Q458").Value = 60
$ws.Range("R458").Value = "Hortaliza"

# Row 459
$ws.Range("A459").Value = 10
$ws.Range("B459").Value = "Vega Modelo de Temuco"
$ws.Range("C459").Value = "La Araucanía"
$ws.Range("D459").Value = 44487
$ws.Range("D459").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E459").Value = 9
$ws.Range("F459").Value = 100112032
$ws.Range("G459").Value = "Zapallo italiano"
$ws.Range("H459").Value = "Sin especificar"
$ws.Range("I459").Value = "Primera"
$ws.Range("J459").Value = 215
$ws.Range("K459").Value = 17000
$ws.Range("L459").Value = 17000
$ws.Range("M459").Value = 17000
$ws.Range("N459").Value = "`$/caja 60 unidades"
$ws.Range("O459").Value = "Región de O'Higgins"
$ws.Range("P459").Value = 283
$ws.Range("Q459").Value = 60
$ws.Range("R459").Value = "Hortaliza"
